# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" (Exhibition) and "全部类型" (All Types) worksheets, matching the
# regenerated output from the data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value to write into column F, shared by both sheets.
$updates = @{
    2  = 13619
    5  = 535
    6  = 497
    9  = 13863
    10 = 14666
    20 = 13
    21 = 1134
    22 = 119
    24 = 5649
    26 = 1047
    27 = 5380
    30 = 212
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $newValue = $updates[$row]
        $cell = $ws.Range("F$row")

        # "全部类型" has an extra row (row 4) compared to "展览", so from row 5
        # downward its matching rows are shifted by one relative to "展览".
        if ($sheetName -eq "全部类型" -and $row -ge 5) {
            $cell = $ws.Range("F" + ($row + 1))
        }

        $cell.Value = $newValue
    }
}
